# "Generate Report for Archive"
#
# The report rows for the two records
#   e0dcfb67-e9cf-4266-acbd-1203e67f0197  (currently row 3)
#   546b8a45-a4fe-43f9-8570-96e9c4393b0d  (currently row 4)
# need to trade places (546b8a45 now sorts/reports before e0dcfb67),
# while the per-row "Status" column and the hyperlink relationship-id
# stay attached to the row position.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Overview": columns A (File Name), B (Path And Name, a
# hyperlink), G (Latest HO Xliff Generate Date) swap between row 3
# and row 4. Columns C/D/E/F keep their current row values.
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$ov_A3 = $wsOverview.Range("A3").Text
$ov_B3 = $wsOverview.Range("B3").Text
$ov_G3 = $wsOverview.Range("G3").Text

$ov_A4 = $wsOverview.Range("A4").Text
$ov_B4 = $wsOverview.Range("B4").Text
$ov_G4 = $wsOverview.Range("G4").Text

$wsOverview.Range("A3").Value = $ov_A4
$wsOverview.Range("B3").Value = $ov_B4
$wsOverview.Range("G3").Value = $ov_G4

$wsOverview.Range("A4").Value = $ov_A3
$wsOverview.Range("B4").Value = $ov_B3
$wsOverview.Range("G4").Value = $ov_G3

foreach ($h in $wsOverview.Hyperlinks) {
    if ($h.Range.Address() -eq '$B$3') {
        $h.TextToDisplay = $ov_B4
    } elseif ($h.Range.Address() -eq '$B$4') {
        $h.TextToDisplay = $ov_B3
    }
}

# ---------------------------------------------------------------
# Sheets "zh-cn" and "de-de": columns A (Source File Name), G
# (Latest Handoff File) and H (Latest Handoff Datetime) swap between
# row 3 and row 4. Column C (Status) keeps its current row value.
# ---------------------------------------------------------------
foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $A3 = $ws.Range("A3").Text
    $G3 = $ws.Range("G3").Text
    $H3 = $ws.Range("H3").Text

    $A4 = $ws.Range("A4").Text
    $G4 = $ws.Range("G4").Text
    $H4 = $ws.Range("H4").Text

    $ws.Range("A3").Value = $A4
    $ws.Range("G3").Value = $G4
    $ws.Range("H3").Value = $H4

    $ws.Range("A4").Value = $A3
    $ws.Range("G4").Value = $G3
    $ws.Range("H4").Value = $H3

    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Address() -eq '$A$3') {
            $h.TextToDisplay = $A4
        } elseif ($h.Range.Address() -eq '$A$4') {
            $h.TextToDisplay = $A3
        }
    }
}
